$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13-17 down to 14-18
$ws.Rows("13:13").Insert()

# Fill in the newly inserted row 13 with the new data record
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44491
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112028
$ws.Range("G13").Value = "Sandia"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 800
$ws.Range("L13").Value = 800
$ws.Range("M13").Value = 800
$ws.Range("N13").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 800
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
